$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("B1").Value = "Colorado"
$ws.Range("C1").Value = 45267
$ws.Range("C1").NumberFormat = "mm-dd-yy"
